$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 160.13333
$ws.Range("I11").Value = 160.13333
$ws.Range("K11").Value = 160.13333
$ws.Range("M11").Value = -20.13333
$ws.Range("H17").Value = 4444.2104
$ws.Range("J17").Value = 4444.2104
$ws.Range("L17").Value = 13332.6312
$ws.Range("N17").Value = -13668.6312
$ws.Range("H58").Value = 803.6667
$ws.Range("J58").Value = 2250
$ws.Range("L58").Value = 6750
$ws.Range("N58").Value = -7050
$ws.Range("H80").Value = 6391.5386
$ws.Range("J80").Value = 7448.6665
$ws.Range("L80").Value = 22345.9995
$ws.Range("N80").Value = -24341.9995
$ws.Range("H83").Value = 6391.5386
$ws.Range("J83").Value = 7448.6665
$ws.Range("L83").Value = 67037.9985
$ws.Range("N83").Value = -77021.9985
$ws.Range("H112").Value = 3537
$ws.Range("J112").Value = 3999.5
$ws.Range("L112").Value = 11998.5
$ws.Range("N112").Value = -14214.5
$ws.Range("H138").Value = 4142.8335
$ws.Range("I138").Value = 2660.2
$ws.Range("K138").Value = 7980.599999999999
$ws.Range("M138").Value = -2840.599999999999

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4965.857
$ws.Range("J63").Value = 4376.5
$ws.Range("L63").Value = 4376.5
$ws.Range("N63").Value = -5748.5
$ws.Range("H66").Value = 4965.857
$ws.Range("J66").Value = 4376.5
$ws.Range("L66").Value = 21882.5
$ws.Range("N66").Value = -28746.5
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2400
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 2400
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H94").Value = 9652
$ws.Range("I94").Value = 9652
$ws.Range("K94").Value = 9652
$ws.Range("M94").Value = -9201

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2780.111
$ws.Range("I58").Value = 2920.1667
$ws.Range("J58").Value = 2500
$ws.Range("K58").Value = 2920.1667
$ws.Range("L58").Value = 2500
$ws.Range("M58").Value = -2717.1667
$ws.Range("N58").Value = -2906
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H99").Value = 2205
$ws.Range("I99").Value = 2012
$ws.Range("K99").Value = 2012
$ws.Range("M99").Value = -514
$ws.Range("H126").Value = 2205
$ws.Range("I126").Value = 2012
$ws.Range("K126").Value = 6036
$ws.Range("M126").Value = -3566
$ws.Range("H136").Value = 2780.111
$ws.Range("I136").Value = 2920.1667
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 8760.500100000001
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -6210.500100000001
$ws.Range("N136").Value = -12600

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 11.777778
$ws.Range("J12").Value = 17
$ws.Range("L12").Value = 51
$ws.Range("N12").Value = -397
$ws.Range("H34").Value = 5121.2
$ws.Range("J34").Value = 8335.333000000001
$ws.Range("L34").Value = 25005.999
$ws.Range("N34").Value = -25173.999
$ws.Range("H68").Value = 6959055.5
$ws.Range("I68").Value = 499
$ws.Range("J68").Value = 7828875
$ws.Range("K68").Value = 1497
$ws.Range("L68").Value = 23486625
$ws.Range("M68").Value = -686
$ws.Range("N68").Value = -23488247
$ws.Range("H71").Value = 6959055.5
$ws.Range("I71").Value = 499
$ws.Range("J71").Value = 7828875
$ws.Range("K71").Value = 4491
$ws.Range("L71").Value = 70459875
$ws.Range("M71").Value = -435
$ws.Range("N71").Value = -70467987
$ws.Range("H92").Value = 234.32259
$ws.Range("I92").Value = 169.90909
$ws.Range("J92").Value = 391.77777
$ws.Range("K92").Value = 509.72727
$ws.Range("L92").Value = 1175.33331
$ws.Range("M92").Value = 738.27273
$ws.Range("N92").Value = -3671.33331
$ws.Range("H132").Value = 3930.5454
$ws.Range("J132").Value = 3319.5
$ws.Range("L132").Value = 29875.5
$ws.Range("N132").Value = -34935.5
$ws.Range("H141").Value = 6067.857
$ws.Range("I141").Value = 6067.857
$ws.Range("K141").Value = 18203.571
$ws.Range("M141").Value = -13023.571

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H97").Value = 249.5
$ws.Range("I97").Value = 269.4
$ws.Range("J97").Value = 150
$ws.Range("K97").Value = 269.4
$ws.Range("L97").Value = 150
$ws.Range("M97").Value = 226.6
$ws.Range("N97").Value = -1142
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 850.8
$ws.Range("I16").Value = 448.42856
$ws.Range("K16").Value = 448.42856
$ws.Range("M16").Value = -278.42856
$ws.Range("H55").Value = 249.95454
$ws.Range("I55").Value = 211.3
$ws.Range("K55").Value = 211.3
$ws.Range("M55").Value = -38.30000000000001
$ws.Range("H128").Value = 21000
$ws.Range("J128").Value = 21000
$ws.Range("L128").Value = 21000
$ws.Range("N128").Value = -30960

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("N28").ClearContents()
$ws.Range("H29").Value = 5333
$ws.Range("I29").Value = 4999.5
$ws.Range("K29").Value = 4999.5
$ws.Range("M29").Value = -4709.5
$ws.Range("H45").Value = 21204
$ws.Range("I45").Value = 21204
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 21204
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -20713
$ws.Range("N45").ClearContents()
$ws.Range("H81").Value = 2000989.4
$ws.Range("I81").Value = 1315
$ws.Range("J81").Value = 5000501
$ws.Range("K81").Value = 2630
$ws.Range("L81").Value = 10001002
$ws.Range("M81").Value = -1569
$ws.Range("N81").Value = -10003124
$ws.Range("H84").Value = 2000989.4
$ws.Range("I84").Value = 1315
$ws.Range("J84").Value = 5000501
$ws.Range("K84").Value = 13150
$ws.Range("L84").Value = 50005010
$ws.Range("M84").Value = -7846
$ws.Range("N84").Value = -50015618
$ws.Range("H119").Value = 119999.5
$ws.Range("J119").Value = 119999.5
$ws.Range("L119").Value = 119999.5
$ws.Range("N119").Value = -129675.5
$ws.Range("H132").Value = 1898.5834
$ws.Range("I132").Value = 523.125
$ws.Range("J132").Value = 4649.5
$ws.Range("K132").Value = 1569.375
$ws.Range("L132").Value = 13948.5
$ws.Range("M132").Value = 960.625
$ws.Range("N132").Value = -19008.5
$ws.Range("H136").Value = 2091.7693
$ws.Range("I136").Value = 2219.6365
$ws.Range("K136").Value = 6658.9095
$ws.Range("M136").Value = -4108.9095
